$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen": update Maximo value for Z3 (C2) ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("C2").Value = "413,1"

# --- Sheet "Metricas": update Tiempo column (B) per zone ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = "412,9"  # Z1: 413,2 -> 412,9
$wsMetricas.Range("B3").Value = "334,0"  # Z2: 334,7 -> 334,0
$wsMetricas.Range("B4").Value = "413,1"  # Z3: 413,3 -> 413,1
$wsMetricas.Range("B5").Value = "336,0"  # Z4: 334,9 -> 336,0

# --- Sheet "Solucion": update reshuffled Pedido/Salida pairs ---
$wsSolucion = $wb.Worksheets.Item("Solucion")
$wsSolucion.Range("B2").Value = "S065"
$wsSolucion.Range("B3").Value = "S025"
$wsSolucion.Range("B7").Value = "S045"
$wsSolucion.Range("B8").Value = "S069"
$wsSolucion.Range("B11").Value = "S042"
$wsSolucion.Range("B13").Value = "S026"
$wsSolucion.Range("B15").Value = "S030"
$wsSolucion.Range("B16").Value = "S070"
$wsSolucion.Range("B17").Value = "S046"
$wsSolucion.Range("B19").Value = "S043"
$wsSolucion.Range("B21").Value = "S067"
$wsSolucion.Range("B22").Value = "S007"
$wsSolucion.Range("B23").Value = "S031"
$wsSolucion.Range("B24").Value = "S071"
$wsSolucion.Range("B25").Value = "S047"
$wsSolucion.Range("B31").Value = "S032"
$wsSolucion.Range("B32").Value = "S048"
$wsSolucion.Range("A36").Value = "Pedido_21"
$wsSolucion.Range("B36").Value = "S049"
$wsSolucion.Range("A37").Value = "Pedido_9"
$wsSolucion.Range("B37").Value = "S033"
$wsSolucion.Range("B38").Value = "S013"
$wsSolucion.Range("B39").Value = "S037"
$wsSolucion.Range("B40").Value = "S053"
$wsSolucion.Range("B41").Value = "S077"
$wsSolucion.Range("B42").Value = "S034"
$wsSolucion.Range("B45").Value = "S010"
$wsSolucion.Range("B46").Value = "S054"
$wsSolucion.Range("B48").Value = "S014"
$wsSolucion.Range("B49").Value = "S038"
$wsSolucion.Range("B50").Value = "S035"
$wsSolucion.Range("B51").Value = "S051"
$wsSolucion.Range("B52").Value = "S075"
$wsSolucion.Range("B53").Value = "S011"
$wsSolucion.Range("B55").Value = "S079"
$wsSolucion.Range("B56").Value = "S039"
$wsSolucion.Range("A57").Value = "Pedido_48"
$wsSolucion.Range("B57").Value = "S055"
$wsSolucion.Range("A58").Value = "Pedido_30"
$wsSolucion.Range("B58").Value = "S036"
$wsSolucion.Range("B59").Value = "S076"
$wsSolucion.Range("B60").Value = "S052"
$wsSolucion.Range("B63").Value = "S040"
$wsSolucion.Range("B64").Value = "S016"
$wsSolucion.Range("B66").Value = "S017"
$wsSolucion.Range("A67").Value = "Pedido_54"
$wsSolucion.Range("B67").Value = "S057"
$wsSolucion.Range("A68").Value = "Pedido_70"
$wsSolucion.Range("B68").Value = "S021"
$wsSolucion.Range("B69").Value = "S061"
$wsSolucion.Range("B70").Value = "S018"
$wsSolucion.Range("B71").Value = "S058"
